# Scheduled-runner price refresh for the Leve profit tracker.
# Updates currentAveragePrice / LevePrice / LeveProfit columns (H:N) on each
# job sheet with freshly pulled market-board figures.
#
# Columns (row 1 header, all sheets):
#   H = currentAveragePrice        I = currentAveragePriceNQ
#   J = currentAveragePriceHQ      K = LevePriceNQ
#   L = LevePriceHQ                M = LeveProfitNQ
#   N = LeveProfitHQ

$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H7").Value = 1865.3334
$ws.Range("J7").Value = 1865.3334
$ws.Range("L7").Value = 1865.3334
$ws.Range("N7").Value = -2089.3334
$ws.Range("H8").Value = 36.42857
$ws.Range("I8").Value = 36.42857
$ws.Range("K8").Value = 109.28571
$ws.Range("M8").Value = 29.71429000000001
$ws.Range("H14").Value = 1865.3334
$ws.Range("J14").Value = 1865.3334
$ws.Range("L14").Value = 1865.3334
$ws.Range("N14").Value = -2247.3334
$ws.Range("H111").Value = 2142.4285
$ws.Range("I111").Value = 2332.8333
$ws.Range("J111").Value = 1000
$ws.Range("K111").Value = 6998.499899999999
$ws.Range("L111").Value = 3000
$ws.Range("M111").Value = -3931.499899999999
$ws.Range("N111").Value = -9134
$ws.Range("H116").Value = 8785.4
$ws.Range("I116").Value = 7475
$ws.Range("J116").Value = 10751
$ws.Range("K116").Value = 7475
$ws.Range("L116").Value = 10751
$ws.Range("M116").Value = -4033
$ws.Range("N116").Value = -17635
$ws.Range("H125").Value = 1416.375
$ws.Range("I125").Value = 2110.3333
$ws.Range("K125").Value = 18992.9997
$ws.Range("M125").Value = -16532.9997
$ws.Range("H132").Value = 1954.3529
$ws.Range("I132").Value = 1873.1428
$ws.Range("K132").Value = 5619.428400000001
$ws.Range("M132").Value = -3089.428400000001
$ws.Range("H137").Value = 3417.8823
$ws.Range("I137").Value = 1996.4445
$ws.Range("K137").Value = 5989.333500000001
$ws.Range("M137").Value = -3439.333500000001
$ws.Range("H138").Value = 3310.6863
$ws.Range("I138").Value = 2708.8572
$ws.Range("K138").Value = 8126.571599999999
$ws.Range("M138").Value = -2986.571599999999

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 18287.416
$ws.Range("I2").Value = 969.4
$ws.Range("K2").Value = 969.4
$ws.Range("M2").Value = -856.4
$ws.Range("H45").Value = 83338104
$ws.Range("I45").Value = 142859150
$ws.Range("K45").Value = 142859150
$ws.Range("M45").Value = -142858773
$ws.Range("H74").Value = 23813210
$ws.Range("I74").Value = 30305132
$ws.Range("J74").Value = 9500
$ws.Range("K74").Value = 30305132
$ws.Range("L74").Value = 9500
$ws.Range("M74").Value = -30304258
$ws.Range("N74").Value = -11248
$ws.Range("H77").Value = 23813210
$ws.Range("I77").Value = 30305132
$ws.Range("J77").Value = 9500
$ws.Range("K77").Value = 151525660
$ws.Range("L77").Value = 47500
$ws.Range("M77").Value = -151521292
$ws.Range("N77").Value = -56236
$ws.Range("H97").Value = 2216.3333
$ws.Range("I97").Value = 2224.75
$ws.Range("J97").Value = 2199.5
$ws.Range("K97").Value = 2224.75
$ws.Range("L97").Value = 2199.5
$ws.Range("M97").Value = -1728.75
$ws.Range("N97").Value = -3191.5
$ws.Range("H110").Value = 5449.081
$ws.Range("I110").Value = 4608
$ws.Range("J110").Value = 7720
$ws.Range("K110").Value = 4608
$ws.Range("L110").Value = 7720
$ws.Range("M110").Value = -2563
$ws.Range("N110").Value = -11810
$ws.Range("H116").Value = 18287.416
$ws.Range("I116").Value = 969.4
$ws.Range("K116").Value = 969.4
$ws.Range("M116").Value = 1324.6
$ws.Range("H132").Value = 3256.8
$ws.Range("I132").Value = 939.43335
$ws.Range("J132").Value = 10208.9
$ws.Range("K132").Value = 2818.30005
$ws.Range("L132").Value = 30626.7
$ws.Range("M132").Value = -288.3000499999998
$ws.Range("N132").Value = -35686.7

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 18287.416
$ws.Range("I3").Value = 969.4
$ws.Range("K3").Value = 969.4
$ws.Range("M3").Value = -855.4
$ws.Range("H20").Value = 6077
$ws.Range("I20").Value = 2745
$ws.Range("J20").Value = 9409
$ws.Range("K20").Value = 2745
$ws.Range("L20").Value = 9409
$ws.Range("M20").Value = -2498
$ws.Range("N20").Value = -9903
$ws.Range("H22").Value = 318.36365
$ws.Range("I22").Value = 322.44446
$ws.Range("K22").Value = 322.44446
$ws.Range("M22").Value = -149.44446
$ws.Range("H107").Value = 5000
$ws.Range("I107").Value = 5000
$ws.Range("K107").Value = 5000
$ws.Range("M107").Value = -3080
$ws.Range("H134").Value = 6590.7
$ws.Range("I134").Value = 1983.3334
$ws.Range("J134").Value = 13501.75
$ws.Range("K134").Value = 5950.0002
$ws.Range("L134").Value = 40505.25
$ws.Range("M134").Value = -3415.0002
$ws.Range("N134").Value = -45575.25

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 2337666
$ws.Range("J4").Value = 2337666
$ws.Range("L4").Value = 2337666
$ws.Range("N4").Value = -2337890
$ws.Range("H7").Value = 371.14285
$ws.Range("I7").Value = 275.25
$ws.Range("K7").Value = 275.25
$ws.Range("M7").Value = -162.25
$ws.Range("H31").Value = 31509.047
$ws.Range("I31").Value = 5316
$ws.Range("K31").Value = 5316
$ws.Range("M31").Value = -5021
$ws.Range("H34").Value = 31509.047
$ws.Range("I34").Value = 5316
$ws.Range("K34").Value = 5316
$ws.Range("M34").Value = -5114
$ws.Range("H132").Value = 5014.25
$ws.Range("I132").Value = 4847.904
$ws.Range("K132").Value = 14543.712
$ws.Range("M132").Value = -12013.712

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H46").Value = 911.7778
$ws.Range("I46").Value = 400.75
$ws.Range("J46").Value = 5000
$ws.Range("K46").Value = 1202.25
$ws.Range("L46").Value = 15000
$ws.Range("M46").Value = -1111.25
$ws.Range("N46").Value = -15182
$ws.Range("H100").Value = 12802.2
$ws.Range("I100").Value = 0
$ws.Range("J100").Value = 12802.2
$ws.Range("K100").Value = 0
$ws.Range("L100").Value = 38406.60000000001
$ws.Range("N100").Value = -40028.60000000001
$ws.Range("M100").ClearContents()
$ws.Range("H112").Value = 125006000
$ws.Range("I112").Value = 166674000
$ws.Range("K112").Value = 500022000
$ws.Range("M112").Value = -500020892
$ws.Range("H117").Value = 2333.0908
$ws.Range("I117").Value = 1011.6667
$ws.Range("J117").Value = 2828.625
$ws.Range("K117").Value = 3035.0001
$ws.Range("L117").Value = 8485.875
$ws.Range("M117").Value = 406.9998999999998
$ws.Range("N117").Value = -15369.875
$ws.Range("H118").Value = 3500
$ws.Range("J118").Value = 6500
$ws.Range("L118").Value = 19500
$ws.Range("N118").Value = -21986
$ws.Range("H119").Value = 20790.545
$ws.Range("I119").Value = 19684.834
$ws.Range("K119").Value = 59054.50199999999
$ws.Range("M119").Value = -54216.50199999999
$ws.Range("H140").Value = 3072.5652
$ws.Range("I140").Value = 2458.2
$ws.Range("K140").Value = 7374.599999999999
$ws.Range("M140").Value = -2194.599999999999

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H93").Value = 35551
$ws.Range("I93").Value = 0
$ws.Range("J93").Value = 35551
$ws.Range("K93").Value = 0
$ws.Range("L93").Value = 35551
$ws.Range("N93").Value = -39295
$ws.Range("M93").ClearContents()

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 10407.765
$ws.Range("I61").Value = 8853.166999999999
$ws.Range("K61").Value = 8853.166999999999
$ws.Range("M61").Value = -8651.166999999999
$ws.Range("H113").Value = 10407.765
$ws.Range("I113").Value = 8853.166999999999
$ws.Range("K113").Value = 8853.166999999999
$ws.Range("M113").Value = -6683.166999999999
$ws.Range("H122").Value = 4713.207
$ws.Range("I122").Value = 3751.12
$ws.Range("J122").Value = 10726.25
$ws.Range("K122").Value = 11253.36
$ws.Range("L122").Value = 32178.75
$ws.Range("M122").Value = -8803.360000000001
$ws.Range("N122").Value = -37078.75
$ws.Range("H130").Value = 63378
$ws.Range("J130").Value = 63378
$ws.Range("L130").Value = 63378
$ws.Range("N130").Value = -73418
$ws.Range("H132").Value = 3105.44
$ws.Range("I132").Value = 1294.1875
$ws.Range("K132").Value = 3882.5625
$ws.Range("M132").Value = -1352.5625

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 379.4
$ws.Range("I107").Value = 266.16666
$ws.Range("K107").Value = 798.4999799999999
$ws.Range("M107").Value = 1121.50002
$ws.Range("H126").Value = 4974.8887
$ws.Range("I126").Value = 4714.4346
$ws.Range("J126").Value = 6472.5
$ws.Range("K126").Value = 14143.3038
$ws.Range("L126").Value = 19417.5
$ws.Range("M126").Value = -11673.3038
$ws.Range("N126").Value = -24357.5
$ws.Range("H132").Value = 4079.1143
$ws.Range("I132").Value = 4292.5557
$ws.Range("K132").Value = 12877.6671
$ws.Range("M132").Value = -10347.6671
$ws.Range("H136").Value = 4886.6665
$ws.Range("I136").Value = 3735.22
$ws.Range("K136").Value = 11205.66
$ws.Range("M136").Value = -8655.66
